$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.2853273333333333
$ws.Range("H2").Value = 0.855982
$ws.Range("I2").Value = 0.05293626506635817
$ws.Range("J2").Value = 0.05293626506635819
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 9.667057999999999
$ws.Range("N2").Value = 29.001174
$ws.Range("O2").Value = 0.1512832311431697
$ws.Range("P2").Value = 0.1588228070066391
$ws.Range("Q2").Value = 2.758275880318666
$ws.Range("R2").Value = 24.824482922868
$ws.Range("S2").Value = 0.008008369223889963
$ws.Range("T2").Value = 0.008407486210286497

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.2853273333333333
$ws.Range("H3").Value = 0.855982
$ws.Range("I3").Value = 0.05293626506635817
$ws.Range("J3").Value = 0.05293626506635819
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 20.60908733333333
$ws.Range("N3").Value = 61.82726199999999
$ws.Range("O3").Value = 0.3225189424433408
$ws.Range("P3").Value = 0.3385924756140875
$ws.Range("Q3").Value = 5.880335931253777
$ws.Range("R3").Value = 52.92302338128399
$ws.Range("S3").Value = 0.0170729482261022
$ws.Range("T3").Value = 0.01792382103858176

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.2853273333333333
$ws.Range("H4").Value = 0.855982
$ws.Range("I4").Value = 0.05293626506635817
$ws.Range("J4").Value = 0.05293626506635819
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 11.51251933333333
$ws.Range("N4").Value = 34.537558
$ws.Range("O4").Value = 0.1801635123472805
$ws.Range("P4").Value = 0.1891424088112641
$ws.Range("Q4").Value = 3.284836441328444
$ws.Range("R4").Value = 29.563527971956
$ws.Range("S4").Value = 0.009537183444901733
$ws.Range("T4").Value = 0.01001249268812256

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.2853273333333333
$ws.Range("H5").Value = 0.855982
$ws.Range("I5").Value = 0.05293626506635817
$ws.Range("J5").Value = 0.05293626506635819
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 13.01136333333333
$ws.Range("N5").Value = 39.03409
$ws.Range("O5").Value = 0.2036194555411202
$ws.Range("P5").Value = 0.2137673372377884
$ws.Range("Q5").Value = 3.712497602931111
$ws.Range("R5").Value = 33.41247842638
$ws.Range("S5").Value = 0.01077885347119227
$ws.Range("T5").Value = 0.01131604442654915

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.2853273333333333
$ws.Range("H6").Value = 0.855982
$ws.Range("I6").Value = 0.05293626506635817
$ws.Range("J6").Value = 0.05293626506635819
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 9.100365500000001
$ws.Range("N6").Value = 18.200731
$ws.Range("O6").Value = 0.1424148585250888
$ws.Range("P6").Value = 0.09967497133022109
$ws.Range("Q6").Value = 2.596583020473667
$ws.Range("R6").Value = 15.579498122842
$ws.Range("S6").Value = 0.007538910700271997
$ws.Range("T6").Value = 0.005276420702818236

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.792877
$ws.Range("H7").Value = 5.378630999999999
$ws.Range("I7").Value = 0.3326292332200106
$ws.Range("J7").Value = 0.3326292332200106
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 9.667057999999999
$ws.Range("N7").Value = 29.001174
$ws.Range("O7").Value = 0.1512832311431697
$ws.Range("P7").Value = 0.1588228070066391
$ws.Range("Q7").Value = 17.331845945866
$ws.Range("R7").Value = 155.986613512794
$ws.Range("S7").Value = 0.05032122517419817
$ws.Range("T7").Value = 0.0528291085124681

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.792877
$ws.Range("H8").Value = 5.378630999999999
$ws.Range("I8").Value = 0.3326292332200106
$ws.Range("J8").Value = 0.3326292332200106
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 20.60908733333333
$ws.Range("N8").Value = 61.82726199999999
$ws.Range("O8").Value = 0.3225189424433408
$ws.Range("P8").Value = 0.3385924756140875
$ws.Range("Q8").Value = 36.94955867092466
$ws.Range("R8").Value = 332.5460280383219
$ws.Range("S8").Value = 0.1072792285238572
$ws.Range("T8").Value = 0.1126257555375791

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.792877
$ws.Range("H9").Value = 5.378630999999999
$ws.Range("I9").Value = 0.3326292332200106
$ws.Range("J9").Value = 0.3326292332200106
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 11.51251933333333
$ws.Range("N9").Value = 34.537558
$ws.Range("O9").Value = 0.1801635123472805
$ws.Range("P9").Value = 0.1891424088112641
$ws.Range("Q9").Value = 20.64053112478866
$ws.Range("R9").Value = 185.764780123098
$ws.Range("S9").Value = 0.05992765096629982
$ws.Range("T9").Value = 0.06291429441227656

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.792877
$ws.Range("H10").Value = 5.378630999999999
$ws.Range("I10").Value = 0.3326292332200106
$ws.Range("J10").Value = 0.3326292332200106
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 13.01136333333333
$ws.Range("N10").Value = 39.03409
$ws.Range("O10").Value = 0.2036194555411202
$ws.Range("P10").Value = 0.2137673372377884
$ws.Range("Q10").Value = 23.32777405897667
$ws.Range("R10").Value = 209.94996653079
$ws.Range("S10").Value = 0.06772978336531883
$ws.Range("T10").Value = 0.07110526547288898

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.792877
$ws.Range("H11").Value = 5.378630999999999
$ws.Range("I11").Value = 0.3326292332200106
$ws.Range("J11").Value = 0.3326292332200106
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 9.100365500000001
$ws.Range("N11").Value = 18.200731
$ws.Range("O11").Value = 0.1424148585250888
$ws.Range("P11").Value = 0.09967497133022109
$ws.Range("Q11").Value = 16.3158359965435
$ws.Range("R11").Value = 97.895015979261
$ws.Range("S11").Value = 0.04737134519033656
$ws.Range("T11").Value = 0.03315480928479798

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 3.311812
$ws.Range("H12").Value = 9.935435999999999
$ws.Range("I12").Value = 0.6144345017136311
$ws.Range("J12").Value = 0.6144345017136311
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 9.667057999999999
$ws.Range("N12").Value = 29.001174
$ws.Range("O12").Value = 0.1512832311431697
$ws.Range("P12").Value = 0.1588228070066391
$ws.Range("Q12").Value = 32.01547868909599
$ws.Range("R12").Value = 288.139308201864
$ws.Range("S12").Value = 0.09295363674508157
$ws.Range("T12").Value = 0.09758621228388449

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 3.311812
$ws.Range("H13").Value = 9.935435999999999
$ws.Range("I13").Value = 0.6144345017136311
$ws.Range("J13").Value = 0.6144345017136311
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 20.60908733333333
$ws.Range("N13").Value = 61.82726199999999
$ws.Range("O13").Value = 0.3225189424433408
$ws.Range("P13").Value = 0.3385924756140875
$ws.Range("Q13").Value = 68.25342273958132
$ws.Range("R13").Value = 614.2808046562319
$ws.Range("S13").Value = 0.1981667656933814
$ws.Range("T13").Value = 0.2080428990379267

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.311812
$ws.Range("H14").Value = 9.935435999999999
$ws.Range("I14").Value = 0.6144345017136311
$ws.Range("J14").Value = 0.6144345017136311
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 11.51251933333333
$ws.Range("N14").Value = 34.537558
$ws.Range("O14").Value = 0.1801635123472805
$ws.Range("P14").Value = 0.1891424088112641
$ws.Range("Q14").Value = 38.12729967836533
$ws.Range("R14").Value = 343.1456971052879
$ws.Range("S14").Value = 0.1106986779360789
$ws.Range("T14").Value = 0.116215621710865

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.311812
$ws.Range("H15").Value = 9.935435999999999
$ws.Range("I15").Value = 0.6144345017136311
$ws.Range("J15").Value = 0.6144345017136311
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 13.01136333333333
$ws.Range("N15").Value = 39.03409
$ws.Range("O15").Value = 0.2036194555411202
$ws.Range("P15").Value = 0.2137673372377884
$ws.Range("Q15").Value = 43.09118922369333
$ws.Range("R15").Value = 387.8207030132399
$ws.Range("S15").Value = 0.125110818704609
$ws.Range("T15").Value = 0.1313460273383503

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.311812
$ws.Range("H16").Value = 9.935435999999999
$ws.Range("I16").Value = 0.6144345017136311
$ws.Range("J16").Value = 0.6144345017136311
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 9.100365500000001
$ws.Range("N16").Value = 18.200731
$ws.Range("O16").Value = 0.1424148585250888
$ws.Range("P16").Value = 0.09967497133022109
$ws.Range("Q16").Value = 30.138699667286
$ws.Range("R16").Value = 180.832198003716
$ws.Range("S16").Value = 0.08750460263448018
$ws.Range("T16").Value = 0.06124374134260486
